# ---------------------------------------------------------------------------
# "Refined metadata to be additional tab"
#
# 1. Refresh the "panel_query_time"-derived F-column timestamps on the
#    existing "data" sheet (rows 2-13) to the new query run.
# 2. Add a new "metadata" worksheet (right after "data") that carries the
#    panel-level metadata (name/id/version/etc.) that used to live only
#    implicitly - now it is its own tab.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- 1. Update the F column timestamps on "data" ---------------------------
$newTimestamps = @(
    "2021-10-05 14:23:06.899064",
    "2021-10-05 14:23:06.899073",
    "2021-10-05 14:23:06.899076",
    "2021-10-05 14:23:06.899079",
    "2021-10-05 14:23:06.899082",
    "2021-10-05 14:23:06.899085",
    "2021-10-05 14:23:06.899088",
    "2021-10-05 14:23:06.899091",
    "2021-10-05 14:23:06.899094",
    "2021-10-05 14:23:06.899097",
    "2021-10-05 14:23:06.899100",
    "2021-10-05 14:23:06.899103"
)

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $newTimestamps[$i]
}

# --- 2. Add the new "metadata" worksheet right after "data" ----------------
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Header row (B1:G1)
$meta.Cells.Item(1, 2).Value = "data_name"
$meta.Cells.Item(1, 3).Value = "data_id"
$meta.Cells.Item(1, 4).Value = "data_version"
$meta.Cells.Item(1, 5).Value = "data_version_created"
$meta.Cells.Item(1, 6).Value = "panel_query_time"
$meta.Cells.Item(1, 7).Value = "panel_get_request"

# Re-use the exact header style/format already used on the "data" tab
# (bold + border + centered) by copying the formatting across.
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

# Data row (A2:G2)
$meta.Cells.Item(2, 1).Value = 0
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$meta.Cells.Item(2, 2).Value = "Vici Syndrome and other autophagy disorders"
$meta.Cells.Item(2, 3).Value = 222

# data_version must stay text ("1.2"), not be coerced to the number 1.2
$dataVersionCell = $meta.Cells.Item(2, 4)
$dataVersionCell.NumberFormat = "@"
$dataVersionCell.Value = "1.2"
$dataVersionCell.Style = "Normal"

$meta.Cells.Item(2, 5).Value = "2019-01-09T13:01:54.387365Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:23:06.895574"
$meta.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/222/?format=json"

# Keep "data" as the active sheet/tab, as it was before this edit.
$ws.Activate()
